$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "758.0 (0.44)"
$ws.Range("E4").Value = "748.0 (0.39)"

$ws.Range("D5").Value = "624.0 (0.56)"
$ws.Range("E5").Value = "611.0 (0.43)"

$ws.Range("D6").Value = "455.0 (0.39)"
$ws.Range("E6").Value = "455.0 (0.34)"

$ws.Range("D7").Value = "920.0 (0.34)"
$ws.Range("E7").Value = "880.0 (0.37)"

$ws.Range("D8").Value = "766.0 (0.26)"
$ws.Range("E8").Value = "702.0 (0.3)"

$ws.Range("D9").Value = "627.0 (0.22)"
$ws.Range("E9").Value = "528.0 (0.26)"

$ws.Range("D10").Value = "757.0 (0.51)"
$ws.Range("E10").Value = "759.0 (0.42)"

$ws.Range("D11").Value = "604.0 (0.52)"
$ws.Range("E11").Value = "612.0 (0.4)"

$ws.Range("D12").Value = "485.0 (1.33)"
$ws.Range("E12").Value = "472.0 (0.82)"
